# "filter short data update"
# Updates the short-data fields (party name, weft/quality/color table, etc.)
# on the active sheet to the new values described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Party name -------------------------------------------------------
$ws.Range("E2").Value = "KAMLESH"

# --- Header block (s no / program date / design) -----------------------
$ws.Range("B4").Value = "1078(3/2)"
$ws.Range("B5").Value = "0000-00-00"
$ws.Range("B6").Value = "M45"

# --- pick / repeat-mtr / panno (numeric-looking text -> keep as text) --
# A leading apostrophe forces these to stay text, matching how the
# original numeric-looking values ("40", "66", ...) were already stored
# as text rather than numbers.
$ws.Range("B7").Value = "'56"
$ws.Range("B10").Value = "'60"

# --- repeat/mtr stays a genuine number ---------------------------------
$ws.Range("B9").Value = 15

# --- weft/deniyar/quality/color rows 5-8 --------------------------------
$ws.Range("F5").Value = "'245"
$ws.Range("G5").Value = "CREKAL"
$ws.Range("H5").Value = "CHIKU"

$ws.Range("F6").Value = "'245"
$ws.Range("G6").Value = "CREKAL"
$ws.Range("H6").Value = "CHIKU"

$ws.Range("F7").Value = "'330"
$ws.Range("G7").Value = "JUTH"
$ws.Range("H7").Value = "VFD9723"

$ws.Range("F8").Value = "'150"
$ws.Range("G8").Value = "MX"
$ws.Range("H8").Value = "PAL MAT"
